$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the data table (before the former 2015 row)
# to hold years 2012-2014, shifting all existing data rows down by 3.
$ws.Rows("2:4").Insert()
$ws.Rows("2:4").ClearFormats()

# Row 2: 2012
$ws.Range("A2").Value = 2012
$ws.Range("B2").Value = 4.27786820933802
$ws.Range("C2").Value = 1.74098131566689
$ws.Range("D2").Value = 8.629802168417978
$ws.Range("E2").Value = 4.263876093965924

# Row 3: 2013
$ws.Range("A3").Value = 2013
$ws.Range("B3").Value = 9.964512723760755
$ws.Range("C3").Value = 3.615055418131874
$ws.Range("D3").Value = 7.030769230769218
$ws.Range("E3").Value = 5.601221522430011

# Row 4: 2014
$ws.Range("A4").Value = 2014
$ws.Range("B4").Value = 5.456550822618689
$ws.Range("C4").Value = -2.054079929300867
$ws.Range("D4").Value = 5.036294379761408
$ws.Range("E4").Value = 1.483826940022026

# Minor value corrections to the (now shifted) 2019 and 2020 rows
$ws.Range("C9").Value = -12.8307694774223
$ws.Range("E9").Value = -0.5450385189066043

$ws.Range("C10").Value = -3.908369428274527
$ws.Range("E10").Value = -2.109051424086783

# Append a new row 14 for year 2024
$ws.Range("A14").Value = 2024
$ws.Range("B14").Value = 10.88236510803766
$ws.Range("C14").Value = -0.43601102235864
$ws.Range("D14").Value = 6.550204742553767
$ws.Range("E14").Value = 5.120555056354159
